$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A8").Value = -21.10860000000001
$ws.Range("A10").Value = -20.47439999999997
$ws.Range("A12").Value = -22.49300000000003
$ws.Range("E13").Value = 12.8387
$ws.Range("A18").Value = -22.32950000000002
$ws.Range("A25").Value = -22.28170000000003
